$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert the new "2022-Q1" fund sheet, positioned right before "总计".
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$ws = $wb.Worksheets.Add($totalSheet)
$ws.Name = "2022-Q1"

# Sheet used purely as a formatting donor (same column layout / styles).
$styleSrc = $wb.Worksheets.Item("2021-Q4")

# --- header row -------------------------------------------------------------
$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

# --- data rows ---------------------------------------------------------------
# Columns B,D,E,F,G hold text-like values (fund codes/ratios kept as strings,
# not numbers) so force Text format before assigning, to avoid Excel
# auto-converting numeric-looking strings (e.g. "005607") into numbers.
$ws.Range("B2:B11").NumberFormat = "@"
$ws.Range("D2:G11").NumberFormat = "@"

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "516970"
$ws.Range("C2").Value = "广发中证基建工程交易型开放式指数证券投资基金"
$ws.Range("D2").Value = "59.01"
$ws.Range("E2").Value = "99.38"
$ws.Range("F2").Value = "2.50"
$ws.Range("G2").Value = "1.4752"
$ws.Range("H2").Value = 10

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "165525"
$ws.Range("C3").Value = "信诚中证基建工程指数（LOF）"
$ws.Range("D3").Value = "17.06"
$ws.Range("E3").Value = "94.00"
$ws.Range("F3").Value = "2.36"
$ws.Range("G3").Value = "0.4026"
$ws.Range("H3").Value = 10

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "005607"
$ws.Range("C4").Value = "华宝中证500指数增强A"
$ws.Range("D4").Value = "0.45"
$ws.Range("E4").Value = "94.72"
$ws.Range("F4").Value = "1.71"
$ws.Range("G4").Value = "0.0077"
$ws.Range("H4").Value = 3

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "005260"
$ws.Range("C5").Value = "银华稳健增利灵活配置混合A"
$ws.Range("D5").Value = "0.32"
$ws.Range("E5").Value = "91.49"
$ws.Range("F5").Value = "1.36"
$ws.Range("G5").Value = "0.0044"
$ws.Range("H5").Value = 2

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "005608"
$ws.Range("C6").Value = "华宝中证500指数增强C"
$ws.Range("D6").Value = "0.23"
$ws.Range("E6").Value = "94.72"
$ws.Range("F6").Value = "1.71"
$ws.Range("G6").Value = "0.0039"
$ws.Range("H6").Value = 3

$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "011987"
$ws.Range("C7").Value = "财通资管智选核心回报6个月持有期混合型发起式证券投资基金A"
$ws.Range("D7").Value = "0.16"
$ws.Range("E7").Value = "38.14"
$ws.Range("F7").Value = "1.52"
$ws.Range("G7").Value = "0.0024"
$ws.Range("H7").Value = 4

$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "009918"
$ws.Range("C8").Value = "上银核心成长混合A"
$ws.Range("D8").Value = "0.13"
$ws.Range("E8").Value = "91.71"
$ws.Range("F8").Value = "0.81"
$ws.Range("G8").Value = "0.0011"
$ws.Range("H8").Value = 6

$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "009919"
$ws.Range("C9").Value = "上银核心成长混合C"
$ws.Range("D9").Value = "0.07"
$ws.Range("E9").Value = "91.71"
$ws.Range("F9").Value = "0.81"
$ws.Range("G9").Value = "0.0006"
$ws.Range("H9").Value = 6

$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "005261"
$ws.Range("C10").Value = "银华稳健增利灵活配置混合C"
$ws.Range("D10").Value = "0.02"
$ws.Range("E10").Value = "91.49"
$ws.Range("F10").Value = "1.36"
$ws.Range("G10").Value = "0.0003"
$ws.Range("H10").Value = 2

$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "011988"
$ws.Range("C11").Value = "财通资管智选核心回报6个月持有期混合型发起式证券投资基金C"
$ws.Range("D11").Value = "0.01"
$ws.Range("E11").Value = "38.14"
$ws.Range("F11").Value = "1.52"
$ws.Range("G11").Value = "0.0002"
$ws.Range("H11").Value = 4

# Column H and A are genuine numbers - already assigned above.

# Strip the incidental "Text" style that NumberFormat="@" attaches, so the
# cells end up with the workbook's default (unstyled) formatting, matching
# the other quarterly sheets.
$ws.Range("B2:B11").Style = "Normal"
$ws.Range("D2:G11").Style = "Normal"

# --- styles: bold/centered/bordered header row + index column ---------------
$styleSrc.Range("B1:H1").Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)

$styleSrc.Range("A2").Copy()
$ws.Range("A2:A11").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) Update the "总计" summary sheet with the new 2022-Q1 totals.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

# Normalize the newly inserted row's formatting from the row below it.
$total.Range("A3:D3").Copy()
$total.Range("A2:D2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 10
$total.Range("D2").Value = 1.9

# Renumber the index column for the rows that shifted down.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5
